{"js": "// Fix abstract indenting: the \"AbstractFirstParagraph\" style carried a\n// direct paragraph-format override (firstLineIndent = 0) that forced the\n// first paragraph of the abstract to have no first-line indent, instead of\n// inheriting the 0.5in (36pt / 720 twips) first-line indent defined on its\n// base style \"Abstract\". Clear that override so the first paragraph of the\n// abstract indents the same way as the rest of the abstract.\nconst styles = context.document.getStyles();\nconst abstractFirstParagraphStyle = styles.getByNameOrNullObject(\"AbstractFirstParagraph\");\nabstractFirstParagraphStyle.load(\"isNullObject\");\nawait context.sync();\n\nif (!abstractFirstParagraphStyle.isNullObject) {\n  const baseStyle = styles.getByNameOrNullObject(\"Abstract\");\n  baseStyle.load(\"isNullObject\");\n  await context.sync();\n\n  // Inherit the first-line indent from the base \"Abstract\" style (falls\n  // back to the documented 36pt/0.5in APA abstract indent if, for some\n  // reason, the base style can't be read).\n  let inheritedFirstLineIndent = 36;\n  if (!baseStyle.isNullObject) {\n    baseStyle.paragraphFormat.load(\"firstLineIndent\");\n    await context.sync();\n    inheritedFirstLineIndent = baseStyle.paragraphFormat.firstLineIndent;\n  }\n\n  abstractFirstParagraphStyle.paragraphFormat.firstLineIndent = inheritedFirstLineIndent;\n  await context.sync();\n}\n", "ps1": "# Fix abstract indenting: the \"AbstractFirstParagraph\" style carried a\n# direct paragraph-format override (FirstLineIndent = 0) that forced the\n# first paragraph of the abstract to have no first-line indent, instead of\n# inheriting the 0.5in (36pt / 720 twips) first-line indent defined on its\n# base style \"Abstract\". Clear that override so the first paragraph of the\n# abstract indents the same way as the rest of the abstract.\n$d = $word.ActiveDocument\n\n$abstractFirstParagraphStyle = $d.Styles(\"AbstractFirstParagraph\")\n\n# Inherit the first-line indent from the base \"Abstract\" style (falls back\n# to the documented 36pt/0.5in APA abstract indent if that lookup fails).\n$inheritedFirstLineIndent = 36\ntry {\n    $baseStyle = $d.Styles(\"Abstract\")\n    $inheritedFirstLineIndent = $baseStyle.ParagraphFormat.FirstLineIndent\n} catch {\n}\n\n$abstractFirstParagraphStyle.ParagraphFormat.FirstLineIndent = $inheritedFirstLineIndent\n"}
